$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Fixed "Date Placeholder" text on the slide master and every slide layout
#    changes from 10/22/2021 -> 10/23/2021 (Insert > Header & Footer > Fixed
#    date, Apply to All).
# ---------------------------------------------------------------------------
$newDate = "10/23/2021"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Reposition the picture ("Grafik 2") on slide 3 - only its vertical
#    offset moves, from y=0 to y=-296427 EMU (x/width/height untouched).
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
    $shp = $slide3.Shapes.Item($i)
    if ($shp.Name -eq "Grafik 2") {
        $shp.Top = -23.34075
    }
}
